# "échange, jour et mois excel"
# Swap columns A ("Mois") and B ("Jour") on the active sheet: the day-of-month
# values (with their centered style) move to column A, and the month names
# move to column B - including the two header cells in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 47

$rngA = $ws.Range("A1:A$lastRow")
$rngB = $ws.Range("B1:B$lastRow")
$scratch = $ws.Range("Z1:Z$lastRow")

# Full copy (values + number formats + styles) via a scratch column so the
# swap is lossless and doesn't synthesize new style indices.
$rngA.Copy($scratch)
$rngB.Copy($rngA)
$scratch.Copy($rngB)
$scratch.ClearContents()

# Re-assert the header cells' own values so the table (Tableau1) picks up
# the new column headers ("Jour" now heads column A, "Mois" heads column B).
$ws.Range("A1").Value = $ws.Range("A1").Value2
$ws.Range("B1").Value = $ws.Range("B1").Value2

# Match the saved selection from the edit.
[void]$ws.Range("E6").Select()
